$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 37452
$ws.Range("D2").Value = 54177121
$ws.Range("C3").Value = 90421
$ws.Range("D3").Value = 132568492
$ws.Range("C4").Value = 30980
$ws.Range("D4").Value = 45883521
$ws.Range("C5").Value = 8643
$ws.Range("D5").Value = 12847735
$ws.Range("C6").Value = 1969
$ws.Range("D6").Value = 2926006
$ws.Range("C11").Value = 40937
$ws.Range("D11").Value = 55568511
$ws.Range("C12").Value = 9583
$ws.Range("D12").Value = 13861146
$ws.Range("C13").Value = 25804
$ws.Range("D13").Value = 37844248
$ws.Range("C14").Value = 8284
$ws.Range("D14").Value = 12295263
$ws.Range("C15").Value = 2138
$ws.Range("D15").Value = 3179383
$ws.Range("C16").Value = 414
$ws.Range("D16").Value = 610123
$ws.Range("C19").Value = 10153
$ws.Range("D19").Value = 13455663
$ws.Range("C20").Value = 13310
$ws.Range("D20").Value = 19223239
$ws.Range("C21").Value = 31508
$ws.Range("D21").Value = 46244898
$ws.Range("C22").Value = 10186
$ws.Range("D22").Value = 15143027
$ws.Range("C23").Value = 2615
$ws.Range("D23").Value = 3890063
$ws.Range("C24").Value = 501
$ws.Range("D24").Value = 745592
$ws.Range("C26").Value = 11609
$ws.Range("D26").Value = 15514208
$ws.Range("C27").Value = 7590
$ws.Range("D27").Value = 10997147
$ws.Range("C28").Value = 22349
$ws.Range("D28").Value = 32803689
$ws.Range("C29").Value = 7766
$ws.Range("D29").Value = 11557302
$ws.Range("C30").Value = 1949
$ws.Range("D30").Value = 2907999
$ws.Range("C31").Value = 362
$ws.Range("D31").Value = 540415
$ws.Range("C33").Value = 8240
$ws.Range("D33").Value = 10891164
$ws.Range("C34").Value = 3196
$ws.Range("D34").Value = 4613190
$ws.Range("C35").Value = 7747
$ws.Range("D35").Value = 11314446
$ws.Range("C36").Value = 3155
$ws.Range("D36").Value = 4675461
$ws.Range("C37").Value = 820
$ws.Range("D37").Value = 1221323
$ws.Range("C38").Value = 161
$ws.Range("D38").Value = 239732
$ws.Range("C40").Value = 2430
$ws.Range("D40").Value = 3285447
$ws.Range("C41").Value = 17110
$ws.Range("D41").Value = 24746829
$ws.Range("C42").Value = 50778
$ws.Range("D42").Value = 74450101
$ws.Range("C43").Value = 18913
$ws.Range("D43").Value = 28094838
$ws.Range("C44").Value = 5578
$ws.Range("D44").Value = 8306978
$ws.Range("C45").Value = 1189
$ws.Range("D45").Value = 1774045
$ws.Range("C49").Value = 16559
$ws.Range("D49").Value = 22062450
$ws.Range("C50").Value = 1978
$ws.Range("D50").Value = 2870316
$ws.Range("C51").Value = 6762
$ws.Range("D51").Value = 9943155
$ws.Range("C52").Value = 2321
$ws.Range("D52").Value = 3466418
$ws.Range("C53").Value = 747
$ws.Range("D53").Value = 1115805
$ws.Range("C56").Value = 6690
$ws.Range("D56").Value = 9217529
$ws.Range("C57").Value = 905
$ws.Range("D57").Value = 1328754
$ws.Range("C58").Value = 2269
$ws.Range("D58").Value = 3367113
$ws.Range("C59").Value = 909
$ws.Range("D59").Value = 1353001
$ws.Range("C60").Value = 311
$ws.Range("D60").Value = 466258
$ws.Range("C63").Value = 1324
$ws.Range("D63").Value = 1869385
$ws.Range("C64").Value = 15243
$ws.Range("D64").Value = 22021293
$ws.Range("C65").Value = 44414
$ws.Range("D65").Value = 65002915
$ws.Range("C66").Value = 15622
$ws.Range("D66").Value = 23220214
$ws.Range("C67").Value = 4546
$ws.Range("D67").Value = 6771292
$ws.Range("C72").Value = 14998
$ws.Range("D72").Value = 19787075
$ws.Range("C73").Value = 50743
$ws.Range("D73").Value = 73850444
$ws.Range("C74").Value = 144489
$ws.Range("D74").Value = 212884868
$ws.Range("C75").Value = 63050
$ws.Range("D75").Value = 93952437
$ws.Range("C76").Value = 20119
$ws.Range("D76").Value = 30059317
$ws.Range("C77").Value = 4749
$ws.Range("D77").Value = 7095223
$ws.Range("C78").Value = 259
$ws.Range("D78").Value = 383670
$ws.Range("C84").Value = 50227
$ws.Range("D84").Value = 68381917
$ws.Range("C85").Value = 4542
$ws.Range("D85").Value = 6579941
$ws.Range("C86").Value = 11448
$ws.Range("D86").Value = 16820065
$ws.Range("C87").Value = 3849
$ws.Range("D87").Value = 5736406
$ws.Range("C88").Value = 1334
$ws.Range("D88").Value = 1993489
$ws.Range("C89").Value = 284
$ws.Range("D89").Value = 423512
$ws.Range("C92").Value = 5319
$ws.Range("D92").Value = 7151162
$ws.Range("C93").Value = 1566
$ws.Range("D93").Value = 2253932
$ws.Range("C94").Value = 5084
$ws.Range("D94").Value = 7490029
$ws.Range("C95").Value = 1926
$ws.Range("D95").Value = 2869446
$ws.Range("C96").Value = 682
$ws.Range("D96").Value = 1021960
$ws.Range("C100").Value = 3478
$ws.Range("D100").Value = 4609562
$ws.Range("C101").Value = 588
$ws.Range("D101").Value = 875664
$ws.Range("C103").Value = 126
$ws.Range("D103").Value = 189000
$ws.Range("C106").Value = 10685
$ws.Range("D106").Value = 15507995
$ws.Range("C107").Value = 29046
$ws.Range("D107").Value = 42680924
$ws.Range("C108").Value = 9733
$ws.Range("D108").Value = 14474213
$ws.Range("C109").Value = 2673
$ws.Range("D109").Value = 3985707
$ws.Range("C113").Value = 9716
$ws.Range("D113").Value = 12845529
$ws.Range("C114").Value = 30173
$ws.Range("D114").Value = 43517343
$ws.Range("C115").Value = 65731
$ws.Range("D115").Value = 96208220
$ws.Range("C116").Value = 21249
$ws.Range("D116").Value = 31579180
$ws.Range("C117").Value = 6024
$ws.Range("D117").Value = 8975826
$ws.Range("C118").Value = 1116
$ws.Range("D118").Value = 1667771
$ws.Range("C123").Value = 25624
$ws.Range("D123").Value = 34241995
$ws.Range("C124").Value = 35633
$ws.Range("D124").Value = 51437821
$ws.Range("C125").Value = 76174
$ws.Range("D125").Value = 111406763
$ws.Range("C126").Value = 23687
$ws.Range("D126").Value = 35157535
$ws.Range("C127").Value = 6338
$ws.Range("D127").Value = 9418551
$ws.Range("C128").Value = 1215
$ws.Range("D128").Value = 1806911
$ws.Range("C132").Value = 31448
$ws.Range("D132").Value = 41786279
$ws.Range("C133").Value = 13145
$ws.Range("D133").Value = 19030046
$ws.Range("C134").Value = 32133
$ws.Range("D134").Value = 47200633
$ws.Range("C135").Value = 11425
$ws.Range("D135").Value = 16975552
$ws.Range("C136").Value = 2942
$ws.Range("D136").Value = 4385805
$ws.Range("C137").Value = 491
$ws.Range("D137").Value = 730490
$ws.Range("C140").Value = 10749
$ws.Range("D140").Value = 14342087
$ws.Range("C141").Value = 34750
$ws.Range("D141").Value = 50193742
$ws.Range("C142").Value = 80643
$ws.Range("D142").Value = 118158264
$ws.Range("C143").Value = 24203
$ws.Range("D143").Value = 35963122
$ws.Range("C144").Value = 6343
$ws.Range("D144").Value = 9465208
$ws.Range("C148").Value = 28927
$ws.Range("D148").Value = 39052638
